$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("TRƯƠNG THÀNH TAM",    "080427", "16:30", "10:34", "B/T",          "RD", "2024-04-11"),
    @("TRƯƠNG TƯ XUÂN",      "080262", "16:30", "10:34", "TEST REQUEST", "RD", "2024-04-11"),
    @("LÊ THANH TUẤN",       "101339", "16:30", "10:37", "TEST REQUEST", "RD", "2024-04-11"),
    @("LÊ PHƯƠNG",           "070032", "16:30", "10:38", "TEST REQUEST", "RD", "2024-04-11"),
    @("TRƯƠNG THÀNH TAM",    "080427", "16:30", "10:38", "B/T",          "RD", "2024-04-11"),
    @("LÊ THANH TUẤN",       "101339", "16:30", "10:39", "TEST REQUEST", "RD", "2024-04-11"),
    @("NGUYỄN HOÀNG VIỆT",   "172759", "16:30", "10:39", "TEST REQUEST", "RD", "2024-04-11"),
    @("BÙI ĐÌNH HỒNG PHÚC",  "193273", "16:30", "10:39", "TEST REQUEST", "RD", "2024-04-11")
)

$startRow = 23
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Columns B (id) and G (date) look numeric/date-like to Excel's
        # auto-detection ("080427", "2024-04-11") and would otherwise be
        # silently coerced to a Number/Date. Force text entry, then restore
        # the default "Normal" style so no stray number-format id is left
        # on the cell (matches the rest of the sheet, which carries no
        # explicit style on data rows).
        if ($c -eq 2 -or $c -eq 7) {
            $cell.NumberFormat = "@"
            $cell.Value = $values[$c - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $values[$c - 1]
        }
    }
}
